$d = $word.ActiveDocument

# --- Edit 1: "Example 1" -> "Example 1 (first stage)" ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Example 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Collapse(0)
    $rng1.InsertAfter(" (first stage)")
}

# --- Edit 2: "Example 2" -> "Example 2 (second stage)" ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Example 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $rng2.InsertAfter(" (second stage)")
}

# --- Edit 3: append the "Example 3" quick-start block + metadata table
#             at the very end of the document body ---
$endRng = $d.Content
$endRng.Collapse(0)
$xmlEnd = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Example </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
        <w:t>3 (third stage)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="TableGrid"/>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="E7E6E6" w:themeColor="background2"/>
          <w:left w:val="single" w:sz="4" w:space="0" w:color="E7E6E6" w:themeColor="background2"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="E7E6E6" w:themeColor="background2"/>
          <w:right w:val="single" w:sz="4" w:space="0" w:color="E7E6E6" w:themeColor="background2"/>
          <w:insideH w:val="single" w:sz="4" w:space="0" w:color="E7E6E6" w:themeColor="background2"/>
          <w:insideV w:val="single" w:sz="4" w:space="0" w:color="E7E6E6" w:themeColor="background2"/>
        </w:tblBorders>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="9350"/>
      </w:tblGrid>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t># Written by June Skeeter</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t># March 2024</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve"># Modified by </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFF00"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>&lt;author&gt;</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve"># Date: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFF00"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>&lt;date&gt;</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>Metadata:</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">  </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>siteID</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFF00"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>SITEID1</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">  </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>estYear</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFF00"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>&lt;YYYY&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">  </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>lat</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFF00"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>&lt;latitude&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">               % North is positive</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">  long: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFF00"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>&lt;longitude&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">             % West is negative</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">  </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>TimeZoneHour</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFF00"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>&lt;</w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFF00"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>timezone</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFF00"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">      % UTC/GMT minus </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>site</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>timezone</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve"> (e.g., for Pacific</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="161"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:ind w:firstLine="3568"/>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>% S</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">tandard </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>T</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve">ime, </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>TimeZoneHour</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve"> = 0 </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>—</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t xml:space="preserve"> 8 = </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>—</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-CA"/>
              </w:rPr>
              <w:t>8)</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>

</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$endRng.InsertXML($xmlEnd)

Write-Host "Edits applied."
